# B6-PowerPoint.pptx edit: 2020-07-27
#
# 1) Three tables (on the slides that contain them) get their
#    <a:tableStyleId> switched from the custom "Table_0" style
#    ({6DCF582B-EA68-4240-AAA2-574A40662F2F}) to the built-in table
#    style {CDF63195-BE9C-499D-AF38-AF353F939070}.
#
# 2) The deck's colour theme is swapped from the "Integral" / "Red
#    Violet" scheme to the standard "Office Theme" palette (this is
#    the theme bound to the slide master / visible slides; it is the
#    half of the theme1.xml <-> theme2.xml swap that is reachable
#    through the PowerPoint object model).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# 1) Re-style every table in the deck.
# ---------------------------------------------------------------
$newTableStyle = "{CDF63195-BE9C-499D-AF38-AF353F939070}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyle)
        }
    }
}

# ---------------------------------------------------------------
# 2) Swap the colour scheme applied to the slide master's theme
#    from "Red Violet" (Integral) to the plain "Office" palette.
# ---------------------------------------------------------------
function HexToRgbVal([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$officeColors = @{
    3  = "44546A"   # dk2
    4  = "E7E6E6"   # lt2
    5  = "5B9BD5"   # accent1
    6  = "ED7D31"   # accent2
    7  = "A5A5A5"   # accent3
    8  = "FFC000"   # accent4
    9  = "4472C4"   # accent5
    10 = "70AD47"   # accent6
    11 = "0563C1"   # hlink
    12 = "954F72"   # folHlink
}

$themeColors = $p.Slides.Item(1).ThemeColorScheme
foreach ($idx in $officeColors.Keys) {
    $themeColors.Item($idx).RGB = HexToRgbVal $officeColors[$idx]
}
